$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1174.7916
$ws.Range("I19").Value = 1125.0834
$ws.Range("K19").Value = 1125.0834
$ws.Range("M19").Value = -950.0834
$ws.Range("H137").Value = 1534.2273
$ws.Range("I137").Value = 928.6429000000001
$ws.Range("K137").Value = 2785.9287
$ws.Range("M137").Value = -235.9287000000004
$ws.Range("H138").Value = 1904.67
$ws.Range("I138").Value = 756.1177
$ws.Range("J138").Value = 2139.9158
$ws.Range("K138").Value = 2268.3531
$ws.Range("L138").Value = 6419.7474
$ws.Range("M138").Value = 2871.6469
$ws.Range("N138").Value = -16699.7474
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 2224.75
$ws.Range("I30").Value = 1449.5
$ws.Range("J30").Value = 3000
$ws.Range("K30").Value = 1449.5
$ws.Range("L30").Value = 3000
$ws.Range("M30").Value = -1299.5
$ws.Range("N30").Value = -3300
$ws.Range("H32").Value = 3562.14
$ws.Range("I32").Value = 3620.2827
$ws.Range("J32").Value = 2893.5
$ws.Range("K32").Value = 3620.2827
$ws.Range("L32").Value = 2893.5
$ws.Range("M32").Value = -3333.2827
$ws.Range("N32").Value = -3467.5
$ws.Range("H45").Value = 1880.9445
$ws.Range("I45").Value = 1833.0714
$ws.Range("J45").Value = 2048.5
$ws.Range("K45").Value = 1833.0714
$ws.Range("L45").Value = 2048.5
$ws.Range("M45").Value = -1456.0714
$ws.Range("N45").Value = -2802.5
$ws.Range("H74").Value = 1376.7333
$ws.Range("I74").Value = 1202.3
$ws.Range("K74").Value = 1202.3
$ws.Range("M74").Value = -328.3
$ws.Range("H77").Value = 1376.7333
$ws.Range("I77").Value = 1202.3
$ws.Range("K77").Value = 6011.5
$ws.Range("M77").Value = -1643.5
$ws.Range("H132").Value = 1787.3636
$ws.Range("I132").Value = 1517.7446
$ws.Range("J132").Value = 3371.375
$ws.Range("K132").Value = 4553.2338
$ws.Range("L132").Value = 10114.125
$ws.Range("M132").Value = -2023.2338
$ws.Range("N132").Value = -15174.125
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4445.5
$ws.Range("I134").Value = 1529.2413
$ws.Range("K134").Value = 4587.7239
$ws.Range("M134").Value = -2052.7239
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1335.9333
$ws.Range("I31").Value = 1041.4166
$ws.Range("J31").Value = 2514
$ws.Range("K31").Value = 1041.4166
$ws.Range("L31").Value = 2514
$ws.Range("M31").Value = -746.4166
$ws.Range("N31").Value = -3104
$ws.Range("H34").Value = 1335.9333
$ws.Range("I34").Value = 1041.4166
$ws.Range("J34").Value = 2514
$ws.Range("K34").Value = 1041.4166
$ws.Range("L34").Value = 2514
$ws.Range("M34").Value = -839.4166
$ws.Range("N34").Value = -2918
$ws.Range("H58").Value = 1101.1613
$ws.Range("I58").Value = 1096.8
$ws.Range("J58").Value = 1109.091
$ws.Range("K58").Value = 1096.8
$ws.Range("L58").Value = 1109.091
$ws.Range("M58").Value = -893.8
$ws.Range("N58").Value = -1515.091
$ws.Range("H132").Value = 6721.52
$ws.Range("I132").Value = 7845.8887
$ws.Range("K132").Value = 23537.6661
$ws.Range("M132").Value = -21007.6661
$ws.Range("H134").Value = 10102396
$ws.Range("I134").Value = 12821838
$ws.Range("J134").Value = 1610.5714
$ws.Range("K134").Value = 38465514
$ws.Range("L134").Value = 4831.7142
$ws.Range("M134").Value = -38462979
$ws.Range("N134").Value = -9901.7142
$ws.Range("H136").Value = 1101.1613
$ws.Range("I136").Value = 1096.8
$ws.Range("J136").Value = 1109.091
$ws.Range("K136").Value = 3290.4
$ws.Range("L136").Value = 3327.273
$ws.Range("M136").Value = -740.3999999999996
$ws.Range("N136").Value = -8427.272999999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 172.625
$ws.Range("I38").Value = 46.833332
$ws.Range("J38").Value = 550
$ws.Range("K38").Value = 140.499996
$ws.Range("L38").Value = 1650
$ws.Range("M38").Value = 206.500004
$ws.Range("N38").Value = -2344
$ws.Range("H68").Value = 1939.4348
$ws.Range("I68").Value = 778.9091
$ws.Range("J68").Value = 3003.25
$ws.Range("K68").Value = 2336.7273
$ws.Range("L68").Value = 9009.75
$ws.Range("M68").Value = -1525.7273
$ws.Range("N68").Value = -10631.75
$ws.Range("H71").Value = 1939.4348
$ws.Range("I71").Value = 778.9091
$ws.Range("J71").Value = 3003.25
$ws.Range("K71").Value = 7010.1819
$ws.Range("L71").Value = 27029.25
$ws.Range("M71").Value = -2954.1819
$ws.Range("N71").Value = -35141.25
$ws.Range("H136").Value = 2196.2727
$ws.Range("I136").Value = 1336
$ws.Range("J136").Value = 2518.875
$ws.Range("K136").Value = 4008
$ws.Range("L136").Value = 7556.625
$ws.Range("M136").Value = 1092
$ws.Range("N136").Value = -17756.625
$ws.Range("H137").Value = 10831
$ws.Range("J137").Value = 14449.643
$ws.Range("L137").Value = 43348.929
$ws.Range("N137").Value = -53548.929
$ws.Range("H138").Value = 2774.3333
$ws.Range("I138").Value = 2721.75
$ws.Range("J138").Value = 2879.5
$ws.Range("K138").Value = 8165.25
$ws.Range("L138").Value = 8638.5
$ws.Range("M138").Value = -3025.25
$ws.Range("N138").Value = -18918.5
$ws.Range("H139").Value = 2588.5
$ws.Range("I139").Value = 3619.875
$ws.Range("J139").Value = 1999.1428
$ws.Range("K139").Value = 10859.625
$ws.Range("L139").Value = 5997.428400000001
$ws.Range("M139").Value = -5719.625
$ws.Range("N139").Value = -16277.4284
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3070.8823
$ws.Range("I80").Value = 1981.8182
$ws.Range("K80").Value = 1981.8182
$ws.Range("M80").Value = -983.8181999999999
$ws.Range("H83").Value = 3070.8823
$ws.Range("I83").Value = 1981.8182
$ws.Range("K83").Value = 9909.091
$ws.Range("M83").Value = -4917.091
$ws.Range("H132").Value = 1698.1951
$ws.Range("I132").Value = 1377.7858
$ws.Range("J132").Value = 2388.3076
$ws.Range("K132").Value = 4133.357400000001
$ws.Range("L132").Value = 7164.9228
$ws.Range("M132").Value = -1603.357400000001
$ws.Range("N132").Value = -12224.9228
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1098.8636
$ws.Range("I16").Value = 887.5
$ws.Range("J16").Value = 2050
$ws.Range("K16").Value = 887.5
$ws.Range("L16").Value = 2050
$ws.Range("M16").Value = -717.5
$ws.Range("N16").Value = -2390
$ws.Range("H22").Value = 1899.6666
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1899.6666
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1899.6666
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -2489.6666
$ws.Range("H27").Value = 1899.6666
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1899.6666
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 1899.6666
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -2113.6666
$ws.Range("H106").Value = 35225
$ws.Range("J106").Value = 35225
$ws.Range("L106").Value = 35225
$ws.Range("N106").Value = -37749
$ws.Range("H136").Value = 6821.278
$ws.Range("I136").Value = 9264.166999999999
$ws.Range("J136").Value = 1935.5
$ws.Range("K136").Value = 27792.501
$ws.Range("L136").Value = 5806.5
$ws.Range("M136").Value = -25242.501
$ws.Range("N136").Value = -10906.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2753.8125
$ws.Range("I81").Value = 2380.75
$ws.Range("J81").Value = 3873
$ws.Range("K81").Value = 4761.5
$ws.Range("L81").Value = 7746
$ws.Range("M81").Value = -3700.5
$ws.Range("N81").Value = -9868
$ws.Range("H84").Value = 2753.8125
$ws.Range("I84").Value = 2380.75
$ws.Range("J84").Value = 3873
$ws.Range("K84").Value = 23807.5
$ws.Range("L84").Value = 38730
$ws.Range("M84").Value = -18503.5
$ws.Range("N84").Value = -49338
$ws.Range("H126").Value = 44445396
$ws.Range("I126").Value = 55556120
$ws.Range("K126").Value = 166668360
$ws.Range("M126").Value = -166665890
$ws.Range("H132").Value = 3860.3572
$ws.Range("I132").Value = 6021
$ws.Range("J132").Value = 1699.7142
$ws.Range("K132").Value = 18063
$ws.Range("L132").Value = 5099.142599999999
$ws.Range("M132").Value = -15533
$ws.Range("N132").Value = -10159.1426
$ws.Range("H136").Value = 582.7778
$ws.Range("I136").Value = 485.8421
$ws.Range("K136").Value = 1457.5263
$ws.Range("M136").Value = 1092.4737
